$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 216, pushing existing rows 216-228 down to 217-229.
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with the new weekly record.
$ws.Range("A216").Value = 8
$ws.Range("B216").Value = "Terminal La Palmera de La Serena"
$ws.Range("C216").Value = "Coquimbo"
$ws.Range("D216").Value = 44516
$ws.Range("E216").Value = 4
$ws.Range("F216").Value = 100114013
$ws.Range("G216").Value = "Zanahoria"
$ws.Range("H216").Value = "Sin especificar"
$ws.Range("I216").Value = "Primera"
$ws.Range("J216").Value = 600
$ws.Range("K216").Value = 6500
$ws.Range("L216").Value = 7000
$ws.Range("M216").Value = 6750
$ws.Range("N216").Value = "`$/saco 20 kilos"
$ws.Range("O216").Value = "Provincia del Elquí"
$ws.Range("P216").Value = 338
$ws.Range("Q216").Value = 20
$ws.Range("R216").Value = "Hortaliza"
